$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values in column D are always stored as plain text in this workbook (e.g.
# "34.776.05" or "230.53"), never as numbers. Cells whose new value would
# otherwise be auto-recognized by Excel as a number are explicitly formatted
# as Text first so the write keeps them as strings, matching the source data.
$ws.Range("D2").Value = "34.776.05"
$ws.Range("D3").Value = "1.825.13"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.53"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.32"
$ws.Range("E8").Value = "  -1.94%  "
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0682"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0987"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").Value = "2.086.03"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.841.48"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.30"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.664"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.64"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "34.605.62"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.35"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "0.0₃0785"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.13"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.13"
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.65"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.24"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.58"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.75"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.123"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.29"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  -8.74%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0548"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.89"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.92"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.23"
$ws.Range("E34").Value = "  +7.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.83"
$ws.Range("E35").Value = "  +2.67%  "
$ws.Range("E36").Value = "  +11.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.697"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "91.49"
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").Value = "1.338.39"
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.02"
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0193"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.45"
$ws.Range("E42").Value = "  -2.04%  "
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.25"
$ws.Range("E44").Value = "  -3.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.75"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.25"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0523"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").Value = "2.001.24"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0670"
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.21"
$ws.Range("E51").Value = "  +13.54%  "
